# Agenda + live demo
#
# 1) Slide 1 (title slide), subtitle shape: split the run
#    "Heidinger, Matthis, " into three runs "Heidinger, " / "Matthis" / ", "
#    (keeps the same look, just separates the first name into its own run,
#    matching how PowerPoint re-flags the word for spell-check).
#
# 2) Slide 2 (Agenda), content placeholder: fix the "Spielstandaes" typo so
#    it reads "... Speicherung des Spielstandes" and move the run boundary
#    so "des" travels with "Spielstandes" instead of with "Speicherung".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 1 - title slide authors
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$subtitle = $slide1.Shapes.Item(2)
$tr1 = $subtitle.TextFrame.TextRange

$full1 = $tr1.Text
$needle1 = "Matthis"
$idx1 = $full1.IndexOf($needle1) + 1
$firstName = $tr1.Characters($idx1, $needle1.Length)
$firstName.Text = $needle1

# ---------------------------------------------------------------------
# Slide 2 - Agenda bullet list
# ---------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$content = $slide2.Shapes.Item(2)
$tr2 = $content.TextFrame.TextRange

$full2 = $tr2.Text
$needle2 = "des Spielstandaes"
$idx2 = $full2.IndexOf($needle2) + 1
$tail = $tr2.Characters($idx2, $needle2.Length)
$tail.Text = "des Spielstandes"
